# Add a new column "Correction " to the Card13 sheet (used as a status/log sheet keyed A1:L13)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card13")

# Header cell M1 - same style as the existing header row (L1) plus the new text
$ws.Range("L1").Copy()
$ws.Range("M1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("M1").Value = "Correction "

# Touch M2:M13 so each row gets a (blank) cell in the new column, matching the
# rest of the sheet which stores an explicit (possibly empty) cell per row.
for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 13).Font.Bold = $false
}

# Row 8 previously had several untouched/blank cells; they now hold the
# placeholder text "nan" like the equivalent cells in every other row.
$ws.Range("D8").Value = "nan"
$ws.Range("F8").Value = "nan"
$ws.Range("G8").Value = "nan"
$ws.Range("H8").Value = "nan"
$ws.Range("J8").Value = "nan"
$ws.Range("K8").Value = "nan"
